# Apply updates described by the commit diff:
#  - Populate O2:Q10 on the "ScoreF" sheet with Minute5/Second5/Rep5 values.
#  - Move the "active"/selected sheet from "ScoreM" to "ScoreF" (tabSelected + activeTab).
#  - Update the cursor/selection on each sheet (ScoreM -> L6, ScoreF -> Q11).

$wb = $excel.ActiveWorkbook

$wsScoreM = $wb.Worksheets.Item("ScoreM")
$wsScoreF = $wb.Worksheets.Item("ScoreF")

# New values for columns O (Minute5), P (Second5), Q (Rep5) on ScoreF, rows 2-10.
$values = @(
    @(8, 0, 208),
    @(8, 0, 121),
    @(8, 0, 107),
    @(8, 0, 163),
    @(8, 0, 179),
    @(8, 0, 152),
    @(8, 0, 197),
    @(8, 0, 248),
    @(8, 0, 239)
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = 2 + $i
    $triplet = $values[$i]
    $wsScoreF.Cells.Item($row, 15).Value = $triplet[0]  # O - Minute5
    $wsScoreF.Cells.Item($row, 16).Value = $triplet[1]  # P - Second5
    $wsScoreF.Cells.Item($row, 17).Value = $triplet[2]  # Q - Rep5
}

# Set the selection (cursor) on ScoreM even though it is no longer the active tab.
$wsScoreM.Range("L6").Select()

# Make ScoreF the active sheet / tab, then set its selection.
$wsScoreF.Activate()
$wsScoreF.Range("Q11").Select()

$wb.Save()
